# "added openpyxl engine support" -- append the five new timeline entries
# (Oct-Dec 2020 measures) that were added to the bottom of the sheet.
#
# Column A holds a date serial (formatted via the sheet's existing "m/d/yyyy"
# style, which is already the default style for column A, so no explicit
# style assignment is required). Column B holds the free-text description,
# which Excel stores as a shared string.
#
# The new shared-string entries must be appended in the same order the
# original author entered them (Gedeeltelijke lockdown, Versoepelingen
# publieke ruimtes, Extra beperkingen op bezoekers, Mondkapjesplicht
# (coronawet), Harde lockdown) -- that order does not match row order
# (row 14 reuses the string that was entered 3rd, row 15 reuses the one
# entered 2nd), so the writes below are intentionally sequenced by
# shared-string insertion order rather than strictly by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 44117
$ws.Range("B13").Value = "Gedeeltelijke lockdown"

$ws.Range("A15").Value = 44152
$ws.Range("B15").Value = "Versoepelingen publieke ruimtes"

$ws.Range("A14").Value = 44138
$ws.Range("B14").Value = "Extra beperkingen op bezoekers"

$ws.Range("A16").Value = 44166
$ws.Range("B16").Value = "Mondkapjesplicht (coronawet)"

$ws.Range("A17").Value = 44179
$ws.Range("B17").Value = "Harde lockdown"

# The author's last selection after typing the new rows was the next empty
# cell in column B.
[void]$ws.Range("B18").Select()
